{"js": "// Fill in the \"Finished\" status cells for the two newly-added PCB-redesign\n// rows (\"No copper zone for L1?\" -> \"YES, NOT DONE YET\" as 3 runs, and\n// \"Maybe move the camera connector to the left\" -> \"YES, DONE\" as 1 run).\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\n// Load the text of the middle (\"Improvement\") column for every row so we can\n// locate the right rows robustly (rather than relying on a fixed index).\nfor (const row of rows.items) {\n  const cells = row.cells.items;\n  const middle = cells[cells.length - 2];\n  middle.load(\"value\");\n}\nawait context.sync();\n\nfunction findRow(snippet) {\n  for (const row of rows.items) {\n    const cells = row.cells.items;\n    const middle = cells[cells.length - 2];\n    if (middle.value && middle.value.indexOf(snippet) !== -1) {\n      return row;\n    }\n  }\n  return null;\n}\n\nconst copperRow = findRow(\"No copper zone for L1\");\nconst cameraRow = findRow(\"move the camera connector to the left\");\n\n// Row \"9\" -> status cell gets three separate runs: \"YES\", \", \", \"NOT DONE YET\".\n// A plain insertText() call merges adjacent inserts into a single run, so\n// build the exact run structure via a flat-OPC insertOoxml() instead.\nif (copperRow) {\n  const statusCell = copperRow.cells.items[copperRow.cells.items.length - 1];\n  const para = statusCell.body.paragraphs.getFirst();\n  const range = para.getRange(\"Start\");\n\n  const flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n        '<pkg:xmlData>' +\n          '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n            '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n          '</Relationships>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n              '<w:p>' +\n                '<w:r><w:t>YES</w:t></w:r>' +\n                '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n                '<w:r><w:t>NOT DONE YET</w:t></w:r>' +\n              '</w:p>' +\n            '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>';\n\n  range.insertOoxml(flatOpc, \"Replace\");\n  await context.sync();\n}\n\n// Row \"10\" -> status cell gets a single run: \"YES, DONE\".\nif (cameraRow) {\n  const statusCell = cameraRow.cells.items[cameraRow.cells.items.length - 1];\n  const para = statusCell.body.paragraphs.getFirst();\n  para.insertText(\"YES, DONE\", \"End\");\n  await context.sync();\n}\n", "ps1": "# Fill in the \"Finished\" status cells for the two newly-added PCB-redesign\n# rows (\"No copper zone for L1?\" -> \"YES, NOT DONE YET\" as 3 runs, and\n# \"Maybe move the camera connector to the left\" -> \"YES, DONE\" as 1 run).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$n = $t.Rows.Count\n\n$copperRow = 0\n$cameraRow = 0\nfor ($i = 1; $i -le $n; $i++) {\n  $middleText = $t.Cell($i, 2).Range.Text\n  if ($middleText -like \"*No copper zone for L1*\") {\n    $copperRow = $i\n  }\n  if ($middleText -like \"*move the camera connector to the left*\") {\n    $cameraRow = $i\n  }\n}\n\n# Row \"9\" -> status cell gets three separate runs: \"YES\", \", \", \"NOT DONE YET\".\n# Assigning plain text to Range.Text merges everything into a single run, so\n# build the exact run structure with a flat-OPC InsertXML() instead.\nif ($copperRow -gt 0) {\n  $cell = $t.Cell($copperRow, 3)\n  $r = $cell.Range\n  $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n        '<pkg:xmlData>' +\n          '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n            '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n          '</Relationships>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n              '<w:p>' +\n                '<w:r><w:t>YES</w:t></w:r>' +\n                '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n                '<w:r><w:t>NOT DONE YET</w:t></w:r>' +\n              '</w:p>' +\n            '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>'\n  $r.InsertXML($xml)\n}\n\n# Row \"10\" -> status cell gets a single run: \"YES, DONE\".\nif ($cameraRow -gt 0) {\n  $cell = $t.Cell($cameraRow, 3)\n  $cell.Range.Text = \"YES, DONE\"\n}\n"}
